$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$qCases = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
WHERE demo.neutered_indicator IN ["Unknown"]  
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$qSamples = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.neutered_indicator IN  ["Unknown"] 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 100
'@

$qFiles = @'
MATCH (f:file)-[*]->(c:case)
MATCH (f)-->(parent)
OPTIONAL MATCH (f)-->(samp:sample)
OPTIONAL MATCH (p:program)<--(s:study)<--(c)
OPTIONAL MATCH (c)<--(demo:demographic)
OPTIONAL MATCH (c)<--(diag:diagnosis)
WITH 
    f, c, parent, samp, p, s, demo, diag
WHERE demo.neutered_indicator IN ["Unknown"] 
WITH 
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    2 as precision
WITH
    DISTINCT f, c, parent, samp, p, s, demo, diag,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH 
    DISTINCT f, c, parent, samp, p, s, demo, diag, unit,
    round(factor * value)/factor AS size
RETURN
    coalesce(f.file_name, '') AS `File Name`,
    coalesce(f.file_format, '') AS `Format`,
    coalesce(f.file_type, '') AS `File Type`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    coalesce(labels(parent)[0], '') AS `Association`,
    coalesce(f.file_description, '') AS `Description`,
    coalesce(samp.sample_id, '') AS `Sample ID`,
    coalesce(c.case_id, '') AS `Case ID`,
    coalesce(demo.breed,'') AS Breed ,
    coalesce(diag.disease_term,'') AS Diagnosis
ORDER BY f.file_name asc
   limit 100
'@

$qStudyFiles = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
WHERE demo.neutered_indicator IN  ["Unknown"] 
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$statQ = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.neutered_indicator IN ["Unknown"]
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Row 1 headers are unchanged: TabName, query, StatQuery, dbExcel, WebExcel

# Tab name for new row
$ws.Range("A5").Value = "StudyFilesTab"

# Neo4j/Web data filenames (D and E columns) - unchanged text, ensure row5 gets them too
$neo4jFile = "TC03_Canine_Filter_NeutStatus-Unknown_Neo4jData.xlsx"
$webFile = "TC03_Canine_Filter_NeutStatus-Unknown_WebData.xlsx"

# Row 2: CasesTab
$ws.Range("B2").Value = $qCases
$ws.Range("C2").Value = $statQ
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

# Row 3: SamplesTab
$ws.Range("B3").Value = $qSamples
$ws.Range("C3").Value = $statQ
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

# Row 4: FilesTab
$ws.Range("B4").Value = $qFiles
$ws.Range("C4").Value = $statQ
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# Row 5: StudyFilesTab (new row)
$ws.Range("B5").Value = $qStudyFiles
$ws.Range("C5").Value = $statQ
$ws.Range("D5").Value = $neo4jFile
$ws.Range("E5").Value = $webFile

# Apply wrap-text style to B and C columns for rows 2-5 (matches existing style used for B2:C4)
$ws.Range("B2:C5").WrapText = $true

# Row heights
$ws.Rows.Item(2).RowHeight = 285
$ws.Rows.Item(3).RowHeight = 238.5
$ws.Rows.Item(4).RowHeight = 235.5
$ws.Rows.Item(5).RowHeight = 234

# Column widths (character units); closest achievable values to target pixel widths.
# Column A is left as auto "best fit" (unchanged) since its content/width did not
# meaningfully change and Excel keeps it as a best-fit column.
$ws.Columns.Item(2).ColumnWidth = 86.66666666666667
$ws.Columns.Item(3).ColumnWidth = 62.5
$ws.Columns.Item(4).ColumnWidth = 69.33333333333333
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668

# Selection / view state
$ws.Range("E5").Select()
$excel.ActiveWindow.Zoom = 60
